$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.637.02"
$ws.Range("E2").Value = "  -1.28%  "

$ws.Range("D3").Value = "1.851.39"
$ws.Range("E3").Value = "  -1.17%  "

$ws.Range("D4").Value = "'1.004"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.21%  "

$ws.Range("D5").Value = "'314.91"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.11%  "

$ws.Range("D6").Value = "'1.004"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.19%  "

$ws.Range("D7").Value = "'0.4257"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.42%  "

$ws.Range("D8").Value = "'0.3650"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.87%  "

$ws.Range("D9").Value = "'44.71"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.60%  "

$ws.Range("D10").Value = "'0.07317"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.38%  "

$ws.Range("D11").Value = "'0.8799"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -6.02%  "

$ws.Range("E12").Value = "  -2.78%  "

$ws.Range("D13").Value = "1.914.99"
$ws.Range("E13").Value = "  -0.21%  "

$ws.Range("D14").Value = "'5.353"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.73%  "

$ws.Range("D15").Value = "'6.544"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.92%  "

$ws.Range("D16").Value = "'0.06932"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.83%  "

$ws.Range("D17").Value = "'1.006"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.23%  "

$ws.Range("D18").Value = "'79.00"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -3.15%  "

$ws.Range("D19").Value = "'0.000008890"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.72%  "

$ws.Range("D20").Value = "'1.004"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.17%  "

$ws.Range("D21").Value = "'15.42"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.00%  "

$ws.Range("D22").Value = "27.644.94"
$ws.Range("E22").Value = "  -1.20%  "

$ws.Range("D23").Value = "'4.999"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.55%  "

$ws.Range("E24").Value = "  -3.68%  "

$ws.Range("D25").Value = "2.096.54"
$ws.Range("E25").Value = "  -1.27%  "

$ws.Range("D26").Value = "'1.984"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.63%  "

$ws.Range("D27").Value = "'153.78"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.76%  "

$ws.Range("D28").Value = "'19.08"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +2.43%  "

$ws.Range("D29").Value = "'122.08"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +7.45%  "

$ws.Range("D30").Value = "'5.258"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -6.26%  "

$ws.Range("D31").Value = "'1.916"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +12.51%  "

$ws.Range("D32").Value = "'0.08946"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.87%  "

$ws.Range("D33").Value = "'0.7620"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -6.26%  "

$ws.Range("D34").Value = "'4.581"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -4.98%  "

$ws.Range("D35").Value = "'2.978"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.25%  "

$ws.Range("D36").Value = "'1.099"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -7.07%  "

$ws.Range("D37").Value = "'1.004"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.21%  "

$ws.Range("B38").Value = "TrustWalletToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D38").Value = "'1.095"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.47%  "

$ws.Range("B39").Value = "Hedera"
$ws.Range("C39").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D39").Value = "'0.05372"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.60%  "

$ws.Range("D40").Value = "'0.01949"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.63%  "

$ws.Range("D41").Value = "'2.811"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -4.96%  "

$ws.Range("D42").Value = "'6.931"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.07%  "

$ws.Range("D43").Value = "'0.5111"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.21%  "

$ws.Range("D44").Value = "'0.1653"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.87%  "

$ws.Range("D45").Value = "'8.281"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -5.79%  "

$ws.Range("D46").Value = "'0.06569"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.73%  "

$ws.Range("D47").Value = "'0.4758"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.81%  "

$ws.Range("D48").Value = "'10.42"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.66%  "

$ws.Range("E49").Value = "  -2.24%  "

$ws.Range("D50").Value = "'1.004"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.20%  "

$ws.Range("E51").Value = "  -2.85%  "
